# Revert file upload functionality
#
# Appends a new row (row 59) at the bottom of each of the four data
# worksheets (FE_LFT_#1, FE_LFT_#2, FE_PLT_#1, FE_PLT_#2), mirroring the
# structure of the existing last row (row 58) on each sheet. Columns B, C,
# E, F, G and I repeat verbatim from row 58 for every sheet, so those are
# copied straight from the existing cells (via .Value2, which preserves the
# exact underlying double for the very large numbers in column G). Columns
# A, D and H carry the new row's own data and are written explicitly.

$wb = $excel.ActiveWorkbook

$newRow = @{
    "FE_LFT_#1" = @{ A = 45845.49159722222; D = "0x01,0x50"; H = 336 }
    "FE_LFT_#2" = @{ A = 45845.49159722222; D = "0x01,0x60"; H = 352 }
    "FE_PLT_#1" = @{ A = 45845.49159722222; D = "0x00,0x68"; H = 104 }
    "FE_PLT_#2" = @{ A = 45845.49159722222; D = "0x00,0x67"; H = 103 }
}

foreach ($sheetName in $newRow.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $newRow[$sheetName]

    $srcRow = 58
    $dstRow = 59

    # Columns identical to the previous last row - copy verbatim.
    $ws.Cells.Item($dstRow, 2).Value = $ws.Range("B$srcRow").Value2
    $ws.Cells.Item($dstRow, 3).Value = $ws.Range("C$srcRow").Value2
    $ws.Cells.Item($dstRow, 5).Value = $ws.Range("E$srcRow").Value2
    $ws.Cells.Item($dstRow, 6).Value = $ws.Range("F$srcRow").Value2
    $ws.Cells.Item($dstRow, 7).Value = $ws.Range("G$srcRow").Value2
    $ws.Cells.Item($dstRow, 9).Value = $ws.Range("I$srcRow").Value2

    # New values for this row.
    $ws.Cells.Item($dstRow, 1).Value = $data.A
    $ws.Cells.Item($dstRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($dstRow, 4).Value = $data.D
    $ws.Cells.Item($dstRow, 8).Value = $data.H
}
